# The presentation's design/theme was switched from the custom "Integral"
# theme to the built-in "Office Theme" (Design > Themes > Office Theme).
# This re-colors the slide master's theme (clrScheme) to the standard
# Office palette while the font/format schemes (already Arial-based in
# both themes) stay the same.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# PpColorSchemeIndex order exposed here: 1=dk1 2=lt1 3=dk2 4=lt2
# 5-10=accent1..accent6 11=hlink 12=folHlink
$officeThemeColors = @(
    0x000000,  # dk1
    0xFFFFFF,  # lt1
    0x44546A,  # dk2
    0xE7E6E6,  # lt2
    0x5B9BD5,  # accent1
    0xED7D31,  # accent2
    0xA5A5A5,  # accent3
    0xFFC000,  # accent4
    0x4472C4,  # accent5
    0x70AD47,  # accent6
    0x0563C1,  # hlink
    0x954F72   # folHlink
)

for ($i = 0; $i -lt $officeThemeColors.Length; $i++) {
    $hex = $officeThemeColors[$i]
    $r = [int]([math]::Floor($hex / 0x10000)) % 0x100
    $g = [int]([math]::Floor($hex / 0x100)) % 0x100
    $b = $hex % 0x100
    $bgr = $r + ($g * 256) + ($b * 65536)
    $scheme.Colors($i + 1).RGB = $bgr
}
